$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''30.934.37'
$ws.Range("E2").Value = '  +3.28%  '

# Row 3
$ws.Range("D3").Value = '''1.911.53'

# Row 4
$ws.Range("E4").Value = '  +0.28%  '

# Row 5
$ws.Range("D5").Value = '''245.92'
$ws.Range("E5").Value = '  +1.04%  '

# Row 6
$ws.Range("E6").Value = '  +0.33%  '

# Row 7
$ws.Range("D7").Value = '''0.4962'
$ws.Range("E7").Value = '  +0.53%  '

# Row 8
$ws.Range("D8").Value = '''0.2995'
$ws.Range("E8").Value = '  +2.94%  '

# Row 9
$ws.Range("D9").Value = '''0.06791'
$ws.Range("E9").Value = '  +2.87%  '

# Row 10
$ws.Range("D10").Value = '''1.913.35'
$ws.Range("E10").Value = '  +1.94%  '

# Row 11
$ws.Range("D11").Value = '''17.04'
$ws.Range("E11").Value = '  +0.69%  '

# Row 12
$ws.Range("D12").Value = '''0.07313'
$ws.Range("E12").Value = '  +1.89%  '

# Row 13
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").Value = '''90.56'
$ws.Range("E13").Value = '  +6.14%  '

# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '''0.6824'
$ws.Range("E14").Value = '  +2.25%  '

# Row 15
$ws.Range("E15").Value = '  +4.87%  '

# Row 16
$ws.Range("D16").Value = '''30.861.45'
$ws.Range("E16").Value = '  +2.93%  '

# Row 17
$ws.Range("D17").Value = '''0.000008019'
$ws.Range("E17").Value = '  +2.22%  '

# Row 18
$ws.Range("D18").Value = '''1.000'
$ws.Range("E18").Value = '  +0.37%  '

# Row 19
$ws.Range("D19").Value = '''13.20'
$ws.Range("E19").Value = '  +3.25%  '

# Row 20
$ws.Range("D20").Value = '''2.158.62'
$ws.Range("E20").Value = '  +1.93%  '

# Row 21
$ws.Range("D21").Value = '''0.9995'
$ws.Range("E21").Value = '  +0.23%  '

# Row 22
$ws.Range("D22").Value = '''4.871'
$ws.Range("E22").Value = '  +2.51%  '

# Row 23
$ws.Range("D23").Value = '''173.48'
$ws.Range("E23").Value = '  +29.13%  '

# Row 24
$ws.Range("D24").Value = '''6.051'
$ws.Range("E24").Value = '  +8.39%  '

# Row 25
$ws.Range("D25").Value = '''9.328'
$ws.Range("E25").Value = '  +2.30%  '

# Row 26
$ws.Range("D26").Value = '''152.52'
$ws.Range("E26").Value = '  +3.31%  '

# Row 27
$ws.Range("D27").Value = '''18.04'
$ws.Range("E27").Value = '  +8.05%  '

# Row 28
$ws.Range("E28").Value = '  +1.44%  '

# Row 29
$ws.Range("E29").Value = '  +2.87%  '

# Row 30
$ws.Range("D30").Value = '''4.326'
$ws.Range("E30").Value = '  +3.64%  '

# Row 31
$ws.Range("D31").Value = '''0.08918'
$ws.Range("E31").Value = '  +3.49%  '

# Row 32
$ws.Range("D32").Value = '''4.075'
$ws.Range("E32").Value = '  +3.84%  '

# Row 33
$ws.Range("D33").Value = '''0.05294'
$ws.Range("E33").Value = '  +5.98%  '

# Row 34
$ws.Range("D34").Value = '''0.7478'
$ws.Range("E34").Value = '  +6.12%  '

# Row 35
$ws.Range("D35").Value = '''1.144'
$ws.Range("E35").Value = '  +3.23%  '

# Row 36
$ws.Range("D36").Value = '''2.644'
$ws.Range("E36").Value = '  -0.26%  '

# Row 37
$ws.Range("D37").Value = '''0.01935'
$ws.Range("E37").Value = '  +17.98%  '

# Row 38
$ws.Range("D38").Value = '''2.725'
$ws.Range("E38").Value = '  +1.38%  '

# Row 39
$ws.Range("D39").Value = '''2.208'
$ws.Range("E39").Value = '  -0.11%  '

# Row 40
$ws.Range("D40").Value = '''0.9410'
$ws.Range("E40").Value = '  +1.05%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.4402'
$ws.Range("E41").Value = '  +5.30%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.975'
$ws.Range("E42").Value = '  -1.36%  '

# Row 43
$ws.Range("D43").Value = '''105.31'

# Row 44
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '''1.001'
$ws.Range("E44").Value = '  +0.72%  '

# Row 45
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '''7.815'
$ws.Range("E45").Value = '  +2.57%  '

# Row 46
$ws.Range("E46").Value = '  +5.56%  '

# Row 47
$ws.Range("D47").Value = '''0.05845'
$ws.Range("E47").Value = '  +2.49%  '

# Row 48
$ws.Range("D48").Value = '''0.3928'
$ws.Range("E48").Value = '  +6.04%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''8.551'
$ws.Range("E49").Value = '  +4.71%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''33.40'
$ws.Range("E50").Value = '  +2.29%  '

# Row 51
$ws.Range("E51").Value = '  +3.36%  '
